# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to H-N columns across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3500  # ALC!H113: 1125 -> 3500
$ws.Cells.Item(113, 9).Value = 3500  # ALC!I113: 1500 -> 3500
$ws.Cells.Item(113, 10).Value = 0  # ALC!J113: 1000 -> 0
$ws.Cells.Item(113, 11).Value = 3500  # ALC!K113: 1500 -> 3500
$ws.Cells.Item(113, 12).Value = 0  # ALC!L113: 1000 -> 0
$ws.Cells.Item(113, 13).Value = -246  # ALC!M113: 1754 -> -246
$ws.Cells.Item(113, 14).ClearContents()  # ALC!N113: -7508 -> (cleared)

$ws.Cells.Item(137, 8).Value = 1375.3903  # ALC!H137: 1357.6428 -> 1375.3903
$ws.Cells.Item(137, 10).Value = 1803.2941  # ALC!J137: 1738.1111 -> 1803.2941
$ws.Cells.Item(137, 12).Value = 5409.8823  # ALC!L137: 5214.3333 -> 5409.8823
$ws.Cells.Item(137, 14).Value = -10509.8823  # ALC!N137: -10314.3333 -> -10509.8823

$ws.Cells.Item(138, 8).Value = 519330.94  # ALC!H138: 495758.6 -> 519330.94
$ws.Cells.Item(138, 9).Value = 1754.6666  # ALC!I138: 1437.0769 -> 1754.6666
$ws.Cells.Item(138, 10).Value = 581440.0600000001  # ALC!J138: 581441 -> 581440.0600000001
$ws.Cells.Item(138, 11).Value = 5263.9998  # ALC!K138: 4311.2307 -> 5263.9998
$ws.Cells.Item(138, 12).Value = 1744320.18  # ALC!L138: 1744323 -> 1744320.18
$ws.Cells.Item(138, 13).Value = -123.9997999999996  # ALC!M138: 828.7692999999999 -> -123.9997999999996
$ws.Cells.Item(138, 14).Value = -1754600.18  # ALC!N138: -1754603 -> -1754600.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101, 8).Value = 35599.8  # ARM!H101: 37000 -> 35599.8
$ws.Cells.Item(101, 10).Value = 35599.8  # ARM!J101: 37000 -> 35599.8
$ws.Cells.Item(101, 12).Value = 35599.8  # ARM!L101: 37000 -> 35599.8
$ws.Cells.Item(101, 14).Value = -42089.8  # ARM!N101: -43490 -> -42089.8

$ws.Cells.Item(110, 8).Value = 1411  # ARM!H110: 1642.3 -> 1411
$ws.Cells.Item(110, 9).Value = 868.55554  # ARM!I110: 982.8333 -> 868.55554
$ws.Cells.Item(110, 11).Value = 868.55554  # ARM!K110: 982.8333 -> 868.55554
$ws.Cells.Item(110, 13).Value = 1176.44446  # ARM!M110: 1062.1667 -> 1176.44446

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 12500688  # BSM!H94: 16667507 -> 12500688
$ws.Cells.Item(94, 9).Value = 20833802  # BSM!I94: 25000516 -> 20833802
$ws.Cells.Item(94, 10).Value = 1017.375  # BSM!J94: 1488 -> 1017.375
$ws.Cells.Item(94, 11).Value = 20833802  # BSM!K94: 25000516 -> 20833802
$ws.Cells.Item(94, 12).Value = 1017.375  # BSM!L94: 1488 -> 1017.375
$ws.Cells.Item(94, 13).Value = -20833351  # BSM!M94: -25000065 -> -20833351
$ws.Cells.Item(94, 14).Value = -1919.375  # BSM!N94: -2390 -> -1919.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(14, 8).Value = 700  # CRP!H14: 0 -> 700
$ws.Cells.Item(14, 10).Value = 700  # CRP!J14: 0 -> 700
$ws.Cells.Item(14, 12).Value = 700  # CRP!L14: 0 -> 700
$ws.Cells.Item(14, 14).Value = -1040  # CRP!N14: None -> -1040

$ws.Cells.Item(15, 8).Value = 1000  # CRP!H15: 1004.5 -> 1000
$ws.Cells.Item(15, 10).Value = 1000  # CRP!J15: 1004.5 -> 1000
$ws.Cells.Item(15, 12).Value = 1000  # CRP!L15: 1004.5 -> 1000
$ws.Cells.Item(15, 14).Value = -1340  # CRP!N15: -1344.5 -> -1340

$ws.Cells.Item(21, 8).Value = 0  # CRP!H21: 3500 -> 0
$ws.Cells.Item(21, 10).Value = 0  # CRP!J21: 3500 -> 0
$ws.Cells.Item(21, 12).Value = 0  # CRP!L21: 3500 -> 0
$ws.Cells.Item(21, 14).ClearContents()  # CRP!N21: -3970 -> (cleared)

$ws.Cells.Item(26, 8).Value = 0  # CRP!H26: 5750 -> 0
$ws.Cells.Item(26, 10).Value = 0  # CRP!J26: 5750 -> 0
$ws.Cells.Item(26, 12).Value = 0  # CRP!L26: 5750 -> 0
$ws.Cells.Item(26, 14).ClearContents()  # CRP!N26: -6324 -> (cleared)

$ws.Cells.Item(32, 8).Value = 1000  # CRP!H32: 1310 -> 1000
$ws.Cells.Item(32, 9).Value = 0  # CRP!I32: 1310 -> 0
$ws.Cells.Item(32, 10).Value = 1000  # CRP!J32: 0 -> 1000
$ws.Cells.Item(32, 11).Value = 0  # CRP!K32: 1310 -> 0
$ws.Cells.Item(32, 12).Value = 1000  # CRP!L32: 0 -> 1000
$ws.Cells.Item(32, 13).ClearContents()  # CRP!M32: -994 -> (cleared)
$ws.Cells.Item(32, 14).Value = -1632  # CRP!N32: None -> -1632

$ws.Cells.Item(33, 8).Value = 5350  # CRP!H33: 0 -> 5350
$ws.Cells.Item(33, 10).Value = 5350  # CRP!J33: 0 -> 5350
$ws.Cells.Item(33, 12).Value = 5350  # CRP!L33: 0 -> 5350
$ws.Cells.Item(33, 14).Value = -6108  # CRP!N33: None -> -6108

$ws.Cells.Item(35, 8).Value = 475  # CRP!H35: 500 -> 475
$ws.Cells.Item(35, 10).Value = 450  # CRP!J35: 0 -> 450
$ws.Cells.Item(35, 12).Value = 450  # CRP!L35: 0 -> 450
$ws.Cells.Item(35, 14).Value = -1038  # CRP!N35: None -> -1038

$ws.Cells.Item(36, 8).Value = 900  # CRP!H36: 500 -> 900
$ws.Cells.Item(36, 9).Value = 0  # CRP!I36: 500 -> 0
$ws.Cells.Item(36, 10).Value = 900  # CRP!J36: 0 -> 900
$ws.Cells.Item(36, 11).Value = 0  # CRP!K36: 500 -> 0
$ws.Cells.Item(36, 12).Value = 900  # CRP!L36: 0 -> 900
$ws.Cells.Item(36, 13).ClearContents()  # CRP!M36: -112 -> (cleared)
$ws.Cells.Item(36, 14).Value = -1676  # CRP!N36: None -> -1676

$ws.Cells.Item(40, 8).Value = 900  # CRP!H40: 500 -> 900
$ws.Cells.Item(40, 9).Value = 0  # CRP!I40: 500 -> 0
$ws.Cells.Item(40, 10).Value = 900  # CRP!J40: 0 -> 900
$ws.Cells.Item(40, 11).Value = 0  # CRP!K40: 500 -> 0
$ws.Cells.Item(40, 12).Value = 900  # CRP!L40: 0 -> 900
$ws.Cells.Item(40, 13).ClearContents()  # CRP!M40: -340 -> (cleared)
$ws.Cells.Item(40, 14).Value = -1220  # CRP!N40: None -> -1220

$ws.Cells.Item(41, 8).Value = 25000  # CRP!H41: 0 -> 25000
$ws.Cells.Item(41, 10).Value = 25000  # CRP!J41: 0 -> 25000
$ws.Cells.Item(41, 12).Value = 25000  # CRP!L41: 0 -> 25000
$ws.Cells.Item(41, 14).Value = -25856  # CRP!N41: None -> -25856

$ws.Cells.Item(42, 8).Value = 5000  # CRP!H42: 4000 -> 5000
$ws.Cells.Item(42, 9).Value = 0  # CRP!I42: 4000 -> 0
$ws.Cells.Item(42, 10).Value = 5000  # CRP!J42: 0 -> 5000
$ws.Cells.Item(42, 11).Value = 0  # CRP!K42: 4000 -> 0
$ws.Cells.Item(42, 12).Value = 5000  # CRP!L42: 0 -> 5000
$ws.Cells.Item(42, 13).ClearContents()  # CRP!M42: -3407 -> (cleared)
$ws.Cells.Item(42, 14).Value = -6186  # CRP!N42: None -> -6186

$ws.Cells.Item(45, 8).Value = 0  # CRP!H45: 3900 -> 0
$ws.Cells.Item(45, 9).Value = 0  # CRP!I45: 3900 -> 0
$ws.Cells.Item(45, 11).Value = 0  # CRP!K45: 3900 -> 0
$ws.Cells.Item(45, 13).ClearContents()  # CRP!M45: -3307 -> (cleared)

$ws.Cells.Item(47, 8).Value = 0  # CRP!H47: 12000 -> 0
$ws.Cells.Item(47, 10).Value = 0  # CRP!J47: 12000 -> 0
$ws.Cells.Item(47, 12).Value = 0  # CRP!L47: 12000 -> 0
$ws.Cells.Item(47, 14).ClearContents()  # CRP!N47: -13132 -> (cleared)

$ws.Cells.Item(50, 8).Value = 27000  # CRP!H50: 20546 -> 27000
$ws.Cells.Item(50, 10).Value = 27000  # CRP!J50: 20546 -> 27000
$ws.Cells.Item(50, 12).Value = 27000  # CRP!L50: 20546 -> 27000
$ws.Cells.Item(50, 14).Value = -28250  # CRP!N50: -21796 -> -28250

$ws.Cells.Item(51, 8).Value = 22155.715  # CRP!H51: 22181.666 -> 22155.715
$ws.Cells.Item(51, 10).Value = 24200  # CRP!J51: 24750 -> 24200
$ws.Cells.Item(51, 12).Value = 24200  # CRP!L51: 24750 -> 24200
$ws.Cells.Item(51, 14).Value = -25672  # CRP!N51: -26222 -> -25672

$ws.Cells.Item(54, 8).Value = 12000  # CRP!H54: 0 -> 12000
$ws.Cells.Item(54, 10).Value = 12000  # CRP!J54: 0 -> 12000
$ws.Cells.Item(54, 12).Value = 12000  # CRP!L54: 0 -> 12000
$ws.Cells.Item(54, 14).Value = -13316  # CRP!N54: None -> -13316

$ws.Cells.Item(55, 8).Value = 3000  # CRP!H55: 8500 -> 3000
$ws.Cells.Item(55, 10).Value = 0  # CRP!J55: 14000 -> 0
$ws.Cells.Item(55, 12).Value = 0  # CRP!L55: 14000 -> 0
$ws.Cells.Item(55, 14).ClearContents()  # CRP!N55: -14630 -> (cleared)

$ws.Cells.Item(57, 8).Value = 19800  # CRP!H57: 0 -> 19800
$ws.Cells.Item(57, 10).Value = 19800  # CRP!J57: 0 -> 19800
$ws.Cells.Item(57, 12).Value = 19800  # CRP!L57: 0 -> 19800
$ws.Cells.Item(57, 14).Value = -20920  # CRP!N57: None -> -20920

$ws.Cells.Item(61, 8).Value = 22155.715  # CRP!H61: 22181.666 -> 22155.715
$ws.Cells.Item(61, 10).Value = 24200  # CRP!J61: 24750 -> 24200
$ws.Cells.Item(61, 12).Value = 24200  # CRP!L61: 24750 -> 24200
$ws.Cells.Item(61, 14).Value = -24896  # CRP!N61: -25446 -> -24896

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 15873919  # CUL!H131: 15385490 -> 15873919
$ws.Cells.Item(131, 9).Value = 58824016  # CUL!I131: 55556044 -> 58824016
$ws.Cells.Item(131, 10).Value = 1056.6086  # CUL!J131: 1023.4894 -> 1056.6086
$ws.Cells.Item(131, 11).Value = 176472048  # CUL!K131: 166668132 -> 176472048
$ws.Cells.Item(131, 12).Value = 3169.8258  # CUL!L131: 3070.4682 -> 3169.8258
$ws.Cells.Item(131, 13).Value = -176467008  # CUL!M131: -166663092 -> -176467008
$ws.Cells.Item(131, 14).Value = -13249.8258  # CUL!N131: -13150.4682 -> -13249.8258

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 19568528  # GSM!H70: 34619016 -> 19568528
$ws.Cells.Item(70, 9).Value = 14709103  # GSM!I70: 35717884 -> 14709103
$ws.Cells.Item(70, 10).Value = 33336900  # GSM!J70: 33337000 -> 33336900
$ws.Cells.Item(70, 11).Value = 14709103  # GSM!K70: 35717884 -> 14709103
$ws.Cells.Item(70, 12).Value = 33336900  # GSM!L70: 33337000 -> 33336900
$ws.Cells.Item(70, 13).Value = -14708833  # GSM!M70: -35717614 -> -14708833
$ws.Cells.Item(70, 14).Value = -33337440  # GSM!N70: -33337540 -> -33337440

$ws.Cells.Item(73, 8).Value = 19568528  # GSM!H73: 34619016 -> 19568528
$ws.Cells.Item(73, 9).Value = 14709103  # GSM!I73: 35717884 -> 14709103
$ws.Cells.Item(73, 10).Value = 33336900  # GSM!J73: 33337000 -> 33336900
$ws.Cells.Item(73, 11).Value = 14709103  # GSM!K73: 35717884 -> 14709103
$ws.Cells.Item(73, 12).Value = 33336900  # GSM!L73: 33337000 -> 33336900
$ws.Cells.Item(73, 13).Value = -14708167  # GSM!M73: -35716948 -> -14708167
$ws.Cells.Item(73, 14).Value = -33338772  # GSM!N73: -33338872 -> -33338772

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2300.6667  # LTW!H7: 2327.5454 -> 2300.6667
$ws.Cells.Item(7, 10).Value = 2321  # LTW!J7: 2400 -> 2321
$ws.Cells.Item(7, 12).Value = 2321  # LTW!L7: 2400 -> 2321
$ws.Cells.Item(7, 14).Value = -2545  # LTW!N7: -2624 -> -2545

$ws.Cells.Item(104, 8).Value = 8246.666999999999  # LTW!H104: 9185 -> 8246.666999999999
$ws.Cells.Item(104, 10).Value = 8246.666999999999  # LTW!J104: 9185 -> 8246.666999999999
$ws.Cells.Item(104, 12).Value = 8246.666999999999  # LTW!L104: 9185 -> 8246.666999999999
$ws.Cells.Item(104, 14).Value = -15234.667  # LTW!N104: -16173 -> -15234.667

$ws.Cells.Item(126, 8).Value = 2300.6667  # LTW!H126: 2327.5454 -> 2300.6667
$ws.Cells.Item(126, 10).Value = 2321  # LTW!J126: 2400 -> 2321
$ws.Cells.Item(126, 12).Value = 6963  # LTW!L126: 7200 -> 6963
$ws.Cells.Item(126, 14).Value = -11903  # LTW!N126: -12140 -> -11903

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 3000  # WVR!H39: 1000 -> 3000
$ws.Cells.Item(39, 9).Value = 0  # WVR!I39: 1000 -> 0
$ws.Cells.Item(39, 10).Value = 3000  # WVR!J39: 0 -> 3000
$ws.Cells.Item(39, 11).Value = 0  # WVR!K39: 1000 -> 0
$ws.Cells.Item(39, 12).Value = 3000  # WVR!L39: 0 -> 3000
$ws.Cells.Item(39, 13).ClearContents()  # WVR!M39: -587 -> (cleared)
$ws.Cells.Item(39, 14).Value = -3826  # WVR!N39: None -> -3826

$ws.Cells.Item(47, 8).Value = 0  # WVR!H47: 12069 -> 0
$ws.Cells.Item(47, 10).Value = 0  # WVR!J47: 12069 -> 0
$ws.Cells.Item(47, 12).Value = 0  # WVR!L47: 12069 -> 0
$ws.Cells.Item(47, 14).ClearContents()  # WVR!N47: -13213 -> (cleared)

$ws.Cells.Item(136, 8).Value = 1514.6666  # WVR!H136: 1559.8096 -> 1514.6666
$ws.Cells.Item(136, 9).Value = 1238.0588  # WVR!I136: 1260.5 -> 1238.0588
$ws.Cells.Item(136, 10).Value = 2186.4285  # WVR!J136: 2158.4285 -> 2186.4285
$ws.Cells.Item(136, 11).Value = 3714.1764  # WVR!K136: 3781.5 -> 3714.1764
$ws.Cells.Item(136, 12).Value = 6559.2855  # WVR!L136: 6475.2855 -> 6559.2855
$ws.Cells.Item(136, 13).Value = -1164.1764  # WVR!M136: -1231.5 -> -1164.1764
$ws.Cells.Item(136, 14).Value = -11659.2855  # WVR!N136: -11575.2855 -> -11659.2855
